# Daily refresh of the "shelf life" tracker sheet.
#
# Columns: A=行号 B=店铺名称 C=地址 D=总天(total days) E=剩余(days remaining)
#          F=开始时间(start date, stored as an integer yyyyMMdd) G/H/I=notes
#
# For every data row (2..last), the "due date" is F + D days. Comparing that
# due date against "today" (2026-01-12, the day this refresh runs):
#   - if the due date has already passed (remaining <= 0), the cycle is
#     restarted: E is reset back to the full D, and F is bumped to today.
#   - otherwise E is simply recomputed as the number of days left.
# Rows whose F value isn't a valid yyyyMMdd date (data-entry typo) are left
# untouched, since there's nothing sane to compute from them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$todayNum = 20260112
$today = [DateTime]::ParseExact([string]$todayNum, "yyyyMMdd", $null)

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {

    $totalDays = [int]$ws.Cells.Item($r, 4).Value()
    $startRaw  = $ws.Cells.Item($r, 6).Value()
    $startStr  = [string]([int]$startRaw)

    if ($startStr -notmatch "^\d{8}$") {
        # malformed start date (e.g. "202510929") -- skip this row entirely
        continue
    }

    $parsedOk = $true
    try {
        $startDate = [DateTime]::ParseExact($startStr, "yyyyMMdd", $null)
    } catch {
        $parsedOk = $false
    }
    if (-not $parsedOk) {
        continue
    }

    $dueDate = $startDate.AddDays($totalDays)
    $remaining = [int]($dueDate.ToOADate() - $today.ToOADate())

    if ($remaining -le 0) {
        $ws.Cells.Item($r, 5).Value = $totalDays
        $ws.Cells.Item($r, 6).Value = $todayNum
    } else {
        $ws.Cells.Item($r, 5).Value = $remaining
    }
}
